$p = $ppt.ActivePresentation

function Find-ShapeById($slide, $targetId) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
    }
    return $null
}

$slide13 = $p.Slides.Item(13)

# --- "Distribution " + "Management " + "System" -> "Distribution " + "Management System"
$mgmtShape = Find-ShapeById $slide13 139
$mgmtRange = $mgmtShape.TextFrame.TextRange
$prefix = "Distribution "
$startPos = $prefix.Length + 1
$tailLen = $mgmtRange.Length - $prefix.Length
$tail = $mgmtRange.Characters($startPos, $tailLen)
$tail.Text = "TEMP_PLACEHOLDER_1"
$mgmtRange2 = $mgmtShape.TextFrame.TextRange
$tail2 = $mgmtRange2.Characters($startPos, "TEMP_PLACEHOLDER_1".Length)
$tail2.Text = "Management System"

# --- "Warehouse" + "Network" -> "WarehouseNetwork"
$whShape = Find-ShapeById $slide13 140
$whRange = $whShape.TextFrame.TextRange
$whRange.Text = "TEMP_PLACEHOLDER_2"
$whShape.TextFrame.TextRange.Text = "WarehouseNetwork"
